# Regenerate sval data to filter save games.
# Updates B,C,D,E,G columns for rows 2-7 on the active sheet with new
# computed values. Column F (Win) is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    3 = @{ B = 0.06328177979961902; C = 0.05231270169004087; D = 0.1529057820181812; E = 0.4998867070740569; G = 0.768386970581898 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    5 = @{ B = 0.7287194209349384; C = 0.3375848360084654;  D = 0.7127328510149897; E = 0.4998867070740569; G = 2.27892381503245 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
